$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new date cells (rows 27-41) pick up the same number format/style
# as the existing Timestamp column before writing values into them.
$ws.Range("B2:B41").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$aValues = @(5509,5420,5377,5365,5329,5315,5301,5299,5237,5235,5218,5225,5309,5257,5343,5300,5262,5271,5258,5357,5451,5528,5604,5658,5857,5932,5960,5948,6026,5993,5982,5982,5928,5849,5862,5802,5688,5634,5597,5557)
$bValues = @(45825,45825.01041666666,45825.02083333334,45825.03125,45825.04166666666,45825.05208333334,45825.0625,45825.07291666666,45825.08333333334,45825.09375,45825.10416666666,45825.11458333334,45825.125,45825.13541666666,45825.14583333334,45825.15625,45825.16666666666,45825.17708333334,45825.1875,45825.19791666666,45825.20833333334,45825.21875,45825.22916666666,45825.23958333334,45825.25,45825.26041666666,45825.27083333334,45825.28125,45825.29166666666,45825.30208333334,45825.3125,45825.32291666666,45825.33333333334,45825.34375,45825.35416666666,45825.36458333334,45825.375,45825.38541666666,45825.39583333334,45825.40625)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
